$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "27.000.85"
Set-TextValue $ws.Range("E2") "  -0.51%  "
Set-TextValue $ws.Range("D3") "1.828.73"
Set-TextValue $ws.Range("E3") "  +0.23%  "
Set-TextValue $ws.Range("D4") "1.006"
Set-TextValue $ws.Range("E4") "  -0.40%  "
Set-TextValue $ws.Range("D5") "312.38"
Set-TextValue $ws.Range("E5") "  -0.02%  "
Set-TextValue $ws.Range("D6") "1.005"
Set-TextValue $ws.Range("E6") "  -0.47%  "
Set-TextValue $ws.Range("D7") "0.4580"
Set-TextValue $ws.Range("E7") "  -0.80%  "
Set-TextValue $ws.Range("D8") "0.3700"
Set-TextValue $ws.Range("E8") "  +1.80%  "
Set-TextValue $ws.Range("D9") "0.07348"
Set-TextValue $ws.Range("E9") "  +0.66%  "
Set-TextValue $ws.Range("D10") "0.8739"
Set-TextValue $ws.Range("E10") "  +0.42%  "
Set-TextValue $ws.Range("B11") "WrappedEther"
Set-TextValue $ws.Range("C11") "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextValue $ws.Range("D11") "1.951.11"
Set-TextValue $ws.Range("E11") "  +6.47%  "
Set-TextValue $ws.Range("B12") "TRON"
Set-TextValue $ws.Range("C12") "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
Set-TextValue $ws.Range("D12") "0.07969"
Set-TextValue $ws.Range("E12") "  +4.96%  "
Set-TextValue $ws.Range("B13") "Solana"
Set-TextValue $ws.Range("C13") "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
Set-TextValue $ws.Range("D13") "19.77"
Set-TextValue $ws.Range("E13") "  -1.71%  "
Set-TextValue $ws.Range("B14") "Polkadot"
Set-TextValue $ws.Range("C14") "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-TextValue $ws.Range("D14") "5.343"
Set-TextValue $ws.Range("E14") "  -0.07%  "
Set-TextValue $ws.Range("B15") "Chainlink"
Set-TextValue $ws.Range("C15") "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
Set-TextValue $ws.Range("D15") "6.544"
Set-TextValue $ws.Range("E15") "  +0.75%  "
Set-TextValue $ws.Range("D16") "91.56"
Set-TextValue $ws.Range("E16") "  -0.98%  "
Set-TextValue $ws.Range("D17") "1.007"
Set-TextValue $ws.Range("E17") "  -0.20%  "
Set-TextValue $ws.Range("D18") "0.000008917"
Set-TextValue $ws.Range("E18") "  +3.23%  "
Set-TextValue $ws.Range("D19") "1.005"
Set-TextValue $ws.Range("E19") "  -0.60%  "
Set-TextValue $ws.Range("E20") "  +2.30%  "
Set-TextValue $ws.Range("D21") "26.870.34"
Set-TextValue $ws.Range("E21") "  -2.06%  "
Set-TextValue $ws.Range("D22") "5.120"
Set-TextValue $ws.Range("E22") "  -1.78%  "
Set-TextValue $ws.Range("B23") "Cosmos"
Set-TextValue $ws.Range("C23") "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextValue $ws.Range("D23") "10.53"
Set-TextValue $ws.Range("E23") "  -0.28%  "
Set-TextValue $ws.Range("B24") "WrappedliquidstakedEther2.0"
Set-TextValue $ws.Range("C24") "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
Set-TextValue $ws.Range("D24") "2.163.20"
Set-TextValue $ws.Range("E24") "  +3.19%  "
Set-TextValue $ws.Range("D25") "153.27"
Set-TextValue $ws.Range("E25") "  +0.95%  "
Set-TextValue $ws.Range("D26") "1.833"
Set-TextValue $ws.Range("E26") "  -2.05%  "
Set-TextValue $ws.Range("D27") "18.40"
Set-TextValue $ws.Range("E27") "  +0.99%  "
Set-TextValue $ws.Range("D28") "2.048"
Set-TextValue $ws.Range("E28") "  -1.51%  "
Set-TextValue $ws.Range("D29") "5.168"
Set-TextValue $ws.Range("E29") "  +1.31%  "
Set-TextValue $ws.Range("D30") "115.48"
Set-TextValue $ws.Range("E30") "  -0.51%  "
Set-TextValue $ws.Range("D31") "0.08880"
Set-TextValue $ws.Range("E31") "  -0.25%  "
Set-TextValue $ws.Range("D32") "2.963"
Set-TextValue $ws.Range("E32") "  +0.04%  "
Set-TextValue $ws.Range("D33") "0.7303"
Set-TextValue $ws.Range("E33") "  -0.26%  "
Set-TextValue $ws.Range("D34") "4.421"
Set-TextValue $ws.Range("E34") "  -0.84%  "
Set-TextValue $ws.Range("E35") "  -0.65%  "
Set-TextValue $ws.Range("D36") "2.462"
Set-TextValue $ws.Range("E36") "  -0.16%  "
Set-TextValue $ws.Range("D37") "1.073"
Set-TextValue $ws.Range("E37") "  -0.24%  "
Set-TextValue $ws.Range("D38") "0.01948"
Set-TextValue $ws.Range("E38") "  +1.78%  "
Set-TextValue $ws.Range("D39") "0.05243"
Set-TextValue $ws.Range("E39") "  -0.05%  "
Set-TextValue $ws.Range("D40") "2.937"
Set-TextValue $ws.Range("E40") "  +0.37%  "
Set-TextValue $ws.Range("D41") "7.119"
Set-TextValue $ws.Range("E41") "  -0.30%  "
Set-TextValue $ws.Range("D42") "0.5146"
Set-TextValue $ws.Range("E42") "  -0.97%  "
Set-TextValue $ws.Range("B43") "Algorand"
Set-TextValue $ws.Range("C43") "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
Set-TextValue $ws.Range("D43") "0.1627"
Set-TextValue $ws.Range("E43") "  -0.09%  "
Set-TextValue $ws.Range("B44") "Aptos"
Set-TextValue $ws.Range("C44") "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextValue $ws.Range("D44") "8.191"
Set-TextValue $ws.Range("E44") "  -0.97%  "
Set-TextValue $ws.Range("B45") "Decentraland"
Set-TextValue $ws.Range("C45") "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
Set-TextValue $ws.Range("D45") "0.4835"
Set-TextValue $ws.Range("E45") "  -0.14%  "
Set-TextValue $ws.Range("B46") "PaxDollar"
Set-TextValue $ws.Range("C46") "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
Set-TextValue $ws.Range("D46") "1.006"
Set-TextValue $ws.Range("E46") "  -0.54%  "
Set-TextValue $ws.Range("B47") "EnergySwap"
Set-TextValue $ws.Range("C47") "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue $ws.Range("D47") "10.15"
Set-TextValue $ws.Range("E47") "  +0.02%  "
Set-TextValue $ws.Range("B48") "Quant"
Set-TextValue $ws.Range("C48") "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
Set-TextValue $ws.Range("D48") "102.44"
Set-TextValue $ws.Range("E48") "  -1.02%  "
Set-TextValue $ws.Range("B49") "NEARProtocol"
Set-TextValue $ws.Range("C49") "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue $ws.Range("D49") "1.634"
Set-TextValue $ws.Range("E49") "  -0.07%  "
Set-TextValue $ws.Range("B50") "Cronos"
Set-TextValue $ws.Range("C50") "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextValue $ws.Range("D50") "0.06200"
Set-TextValue $ws.Range("E50") "  -1.06%  "
Set-TextValue $ws.Range("B51") "Aave"
Set-TextValue $ws.Range("C51") "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextValue $ws.Range("D51") "64.61"
Set-TextValue $ws.Range("E51") "  +0.19%  "
